# Applies the "Batterywise analysis" label/value updates to the
# "Analysis Results" sheet: relabels several metrics with their units,
# fixes a couple of swapped/incorrect values, fills in a previously
# blank cell, shifts the speed-bucket rows down by one, and appends a
# new "Time spent in 80-90 km/h" row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6-7: Starting / Ending SoC values were swapped ---
$ws.Range("B6").Value = 100
$ws.Range("B7").Value = 10

# --- Row 8: relabel distance metric ---
$ws.Range("A8").Value = "Total distance covered (km)"

# --- Row 9: relabel WH/KM metric ---
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# --- Row 10: relabel + correct SOC consumed ---
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 90

# --- Row 12-14: append units to labels ---
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# --- Row 15: relabel + fix sign of Regenerative Effectiveness ---
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 4.408442272039773

# --- Rows 16-17: Highest/Lowest Cell Voltage labels+values were swapped ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.382
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.125

# --- Row 18: append unit ---
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# --- Rows 19-21: append units, fill in missing Difference in Temperature ---
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 18

# --- Rows 22-27: append units to labels ---
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Rows 28-29: lowest/highest cell temp labels were swapped ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# --- Row 30: append unit ---
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: was "Maximum BMS Temperature in C" -> now Battery Voltage(V) ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 55

# --- Row 32: was "Battery Voltage" -> now Total energy charged(kWh) ---
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.830360002777778

# --- Row 33: was "Total energy charged in kWh" -> now Electricity consumption units(kW) ---
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.00000007400776333405215

# --- Row 34: was "Electricity consumption units in kW" -> now Idling time percentage ---
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 2.723747980613894

# --- Row 35: was "Idling time percentage" -> now Time spent in 0-10 km/h ---
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 18.17770597738288

# --- Row 36: was "Time spent in 0-10 km/h" -> now Time spent in 10-20 km/h ---
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.647819063004846

# --- Row 37: was "Time spent in 10-20 km/h" -> now Time spent in 20-30 km/h ---
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 7.846526655896607

# --- Row 38: was "Time spent in 20-30 km/h" -> now Time spent in 30-40 km/h ---
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 45.86914378029079

# --- Row 39: was "Time spent in 30-40 km/h" -> now Time spent in 40-50 km/h ---
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 20.29886914378029

# --- Row 40: was "Time spent in 40-50 km/h" -> now Time spent in 50-60 km/h ---
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 0

# --- Row 41: was "Time spent in 50-60 km/h" -> now Time spent in 60-70 km/h ---
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 0

# --- Row 42: was "Time spent in 60-70 km/h" -> now Time spent in 70-80 km/h ---
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

# --- Row 43: new row appended ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
